# Add GCD, LCM, Integer Division, and Factorial (#30)
# Update the Status column ("D") on the Functions sheet from "N/A" to "Done"
# for the functions that have now been implemented.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Rows whose Status (column D) moves from N/A -> Done
$doneRows = @(6, 8, 9, 11, 12, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43)

foreach ($r in $doneRows) {
    $ws.Range("D$r").Value = "Done"
}

# Move the active selection to D7 (no more scrolled-down view)
$ws.Activate()
$ws.Range("D7").Select() | Out-Null
